$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H19").Value = 645.4545000000001
$ws.Range("I19").Value = 630
$ws.Range("J19").Value = 658.3333
$ws.Range("K19").Value = 630
$ws.Range("L19").Value = 658.3333
$ws.Range("M19").Value = -455
$ws.Range("N19").Value = -1008.3333
$ws.Range("H64").Value = 3411.1135
$ws.Range("I64").Value = 3239.087
$ws.Range("J64").Value = 3599.524
$ws.Range("K64").Value = 3239.087
$ws.Range("L64").Value = 3599.524
$ws.Range("M64").Value = -2991.087
$ws.Range("N64").Value = -4095.524
$ws.Range("H67").Value = 3411.1135
$ws.Range("I67").Value = 3239.087
$ws.Range("J67").Value = 3599.524
$ws.Range("K67").Value = 3239.087
$ws.Range("L67").Value = 3599.524
$ws.Range("M67").Value = -2381.087
$ws.Range("N67").Value = -5315.523999999999
$ws.Range("H130").Value = 49998.332
$ws.Range("J130").Value = 49998.332
$ws.Range("L130").Value = 49998.332
$ws.Range("N130").Value = -60038.332
$ws.Range("H138").Value = 1730.7234
$ws.Range("I138").Value = 1018.27026
$ws.Range("J138").Value = 4366.8
$ws.Range("K138").Value = 3054.81078
$ws.Range("L138").Value = 13100.4
$ws.Range("M138").Value = 2085.18922
$ws.Range("N138").Value = -23380.4

$ws = $wb.Worksheets("ARM")
$ws.Range("H74").Value = 2100
$ws.Range("I74").Value = 1250
$ws.Range("J74").Value = 2666.6667
$ws.Range("K74").Value = 1250
$ws.Range("L74").Value = 2666.6667
$ws.Range("M74").Value = -376
$ws.Range("N74").Value = -4414.6667
$ws.Range("H77").Value = 2100
$ws.Range("I77").Value = 1250
$ws.Range("J77").Value = 2666.6667
$ws.Range("K77").Value = 6250
$ws.Range("L77").Value = 13333.3335
$ws.Range("M77").Value = -1882
$ws.Range("N77").Value = -22069.3335
$ws.Range("H109").Value = 33000
$ws.Range("J109").Value = 33000
$ws.Range("L109").Value = 33000
$ws.Range("N109").Value = -35774
$ws.Range("H122").Value = 2475.05
$ws.Range("I122").Value = 1862.4546
$ws.Range("K122").Value = 5587.3638
$ws.Range("M122").Value = -3137.3638

$ws = $wb.Worksheets("BSM")
$ws.Range("H29").Value = 5398
$ws.Range("I29").Value = 796
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 796
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = -507
$ws.Range("N29").Value = -10578
$ws.Range("H36").Value = 465
$ws.Range("I36").Value = 465
$ws.Range("K36").Value = 465
$ws.Range("M36").Value = 69
$ws.Range("H75").Value = 258654.5
$ws.Range("I75").Value = 42333.332
$ws.Range("J75").Value = 388447.2
$ws.Range("K75").Value = 42333.332
$ws.Range("L75").Value = 388447.2
$ws.Range("M75").Value = -41397.332
$ws.Range("N75").Value = -390319.2
$ws.Range("H78").Value = 258654.5
$ws.Range("I78").Value = 42333.332
$ws.Range("J78").Value = 388447.2
$ws.Range("K78").Value = 126999.996
$ws.Range("L78").Value = 1165341.6
$ws.Range("M78").Value = -122319.996
$ws.Range("N78").Value = -1174701.6
$ws.Range("H139").Value = 129206.664
$ws.Range("I139").Value = 78000
$ws.Range("J139").Value = 139448
$ws.Range("K139").Value = 78000
$ws.Range("L139").Value = 139448
$ws.Range("M139").Value = -72860
$ws.Range("N139").Value = -149728

$ws = $wb.Worksheets("CRP")
$ws.Range("H132").Value = 903770.8
$ws.Range("I132").Value = 1353484.9
$ws.Range("J132").Value = 4342.6
$ws.Range("K132").Value = 4060454.7
$ws.Range("L132").Value = 13027.8
$ws.Range("M132").Value = -4057924.7
$ws.Range("N132").Value = -18087.8

$ws = $wb.Worksheets("CUL")
$ws.Range("H8").Value = 1290.091
$ws.Range("I8").Value = 1290.091
$ws.Range("K8").Value = 3870.273
$ws.Range("M8").Value = -3731.273
$ws.Range("H134").Value = 3551
$ws.Range("I134").Value = 1960
$ws.Range("J134").Value = 6534.125
$ws.Range("K134").Value = 5880
$ws.Range("L134").Value = 19602.375
$ws.Range("M134").Value = -810
$ws.Range("N134").Value = -29742.375
$ws.Range("H138").Value = 3025.375
$ws.Range("I138").Value = 1665
$ws.Range("J138").Value = 3478.8333
$ws.Range("K138").Value = 4995
$ws.Range("L138").Value = 10436.4999
$ws.Range("M138").Value = 145
$ws.Range("N138").Value = -20716.4999

$ws = $wb.Worksheets("GSM")
$ws.Range("H97").Value = 42069.31
$ws.Range("I97").Value = 86998.336
$ws.Range("J97").Value = 3558.7144
$ws.Range("K97").Value = 86998.336
$ws.Range("L97").Value = 3558.7144
$ws.Range("M97").Value = -86502.336
$ws.Range("N97").Value = -4550.7144
$ws.Range("H123").Value = 25478.309
$ws.Range("J123").Value = 25478.309
$ws.Range("L123").Value = 25478.309
$ws.Range("N123").Value = -30378.309
$ws.Range("H131").Value = 48657
$ws.Range("J131").Value = 48657
$ws.Range("L131").Value = 48657
$ws.Range("N131").Value = -58737

$ws = $wb.Worksheets("LTW")
$ws.Range("H61").Value = 23570.7
$ws.Range("J61").Value = 3899.6667
$ws.Range("L61").Value = 3899.6667
$ws.Range("N61").Value = -4303.6667
$ws.Range("H100").Value = 2160.6
$ws.Range("I100").Value = 1767.6666
$ws.Range("J100").Value = 2750
$ws.Range("K100").Value = 1767.6666
$ws.Range("L100").Value = 2750
$ws.Range("M100").Value = -1226.6666
$ws.Range("N100").Value = -3832
$ws.Range("H113").Value = 23570.7
$ws.Range("J113").Value = 3899.6667
$ws.Range("L113").Value = 3899.6667
$ws.Range("N113").Value = -8239.6667
$ws.Range("H122").Value = 75002090
$ws.Range("I122").Value = 125001000
$ws.Range("J122").Value = 50002628
$ws.Range("K122").Value = 375003000
$ws.Range("L122").Value = 150007884
$ws.Range("M122").Value = -375000550
$ws.Range("N122").Value = -150012784

$ws = $wb.Worksheets("WVR")
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 4250
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 4250
$ws.Range("M62").Value = -2876
$ws.Range("N62").Value = -5498
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 4250
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 21250
$ws.Range("M65").Value = -14380
$ws.Range("N65").Value = -27490
$ws.Range("H68").Value = 20000
$ws.Range("J68").Value = 20000
$ws.Range("L68").Value = 20000
$ws.Range("N68").Value = -21622
$ws.Range("H71").Value = 20000
$ws.Range("J71").Value = 20000
$ws.Range("L71").Value = 60000
$ws.Range("N71").Value = -68112
$ws.Range("H96").Value = 1000
$ws.Range("I96").Value = 1000
$ws.Range("K96").Value = 1000
$ws.Range("M96").Value = 373
$ws.Range("H122").Value = 10419817
$ws.Range("J122").Value = 4017.8333
$ws.Range("L122").Value = 12053.4999
$ws.Range("N122").Value = -16953.4999
$ws.Range("H123").Value = 34246
$ws.Range("J123").Value = 34246
$ws.Range("L123").Value = 34246
$ws.Range("N123").Value = -44046
